$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1907.1
$ws.Range("I2").Value = 992.5714
$ws.Range("K2").Value = 992.5714
$ws.Range("M2").Value = -879.5714
$ws.Range("H9").Value = 1949.8334
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 2000
$ws.Range("N40").Value = -2350
$ws.Range("H43").Value = 6682.1904
$ws.Range("I43").Value = 6768.077
$ws.Range("J43").Value = 6542.625
$ws.Range("K43").Value = 6768.077
$ws.Range("L43").Value = 6542.625
$ws.Range("M43").Value = -6699.077
$ws.Range("N43").Value = -6680.625
$ws.Range("H58").Value = 2891.7
$ws.Range("I58").Value = 1071.4286
$ws.Range("J58").Value = 7139
$ws.Range("K58").Value = 3214.2858
$ws.Range("L58").Value = 21417
$ws.Range("M58").Value = -3064.2858
$ws.Range("N58").Value = -21717
$ws.Range("H69").Value = 6000
$ws.Range("J69").Value = 6000
$ws.Range("L69").Value = 18000
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 6000
$ws.Range("J72").Value = 6000
$ws.Range("L72").Value = 54000
$ws.Range("N72").Value = -62736
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
$ws.Range("H88").Value = 2188.0557
$ws.Range("I88").Value = 2390.5
$ws.Range("J88").Value = 2162.75
$ws.Range("K88").Value = 2390.5
$ws.Range("L88").Value = 2162.75
$ws.Range("M88").Value = -1984.5
$ws.Range("N88").Value = -2974.75
$ws.Range("H91").Value = 2188.0557
$ws.Range("I91").Value = 2390.5
$ws.Range("J91").Value = 2162.75
$ws.Range("K91").Value = 2390.5
$ws.Range("L91").Value = 2162.75
$ws.Range("M91").Value = -986.5
$ws.Range("N91").Value = -4970.75
$ws.Range("H111").Value = 1895.3
$ws.Range("I111").Value = 1703
$ws.Range("K111").Value = 5109
$ws.Range("M111").Value = -2042
$ws.Range("H116").Value = 3083.5715
$ws.Range("I116").Value = 2850
$ws.Range("K116").Value = 2850
$ws.Range("M116").Value = 592
$ws.Range("H135").Value = 616.63635
$ws.Range("I135").Value = 316
$ws.Range("J135").Value = 1142.75
$ws.Range("K135").Value = 2844
$ws.Range("L135").Value = 10284.75
$ws.Range("M135").Value = -309
$ws.Range("N135").Value = -15354.75
$ws.Range("H137").Value = 1748
$ws.Range("I137").Value = 1663.4166
$ws.Range("J137").Value = 1874.875
$ws.Range("K137").Value = 4990.2498
$ws.Range("L137").Value = 5624.625
$ws.Range("M137").Value = -2440.2498
$ws.Range("N137").Value = -10724.625
$ws.Range("H138").Value = 1859.6666
$ws.Range("J138").Value = 4000
$ws.Range("L138").Value = 12000
$ws.Range("N138").Value = -22280

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3808.0715
$ws.Range("I32").Value = 2689.5217
$ws.Range("K32").Value = 2689.5217
$ws.Range("M32").Value = -2402.5217
$ws.Range("H61").Value = 798.75
$ws.Range("I61").Value = 798.75
$ws.Range("K61").Value = 798.75
$ws.Range("M61").Value = -586.75
$ws.Range("H97").Value = 1566.4
$ws.Range("I97").Value = 950.5714
$ws.Range("J97").Value = 3003.3333
$ws.Range("K97").Value = 950.5714
$ws.Range("L97").Value = 3003.3333
$ws.Range("M97").Value = -454.5714
$ws.Range("N97").Value = -3995.3333
$ws.Range("H122").Value = 3002.3
$ws.Range("I122").Value = 2257.6667
$ws.Range("K122").Value = 6773.000100000001
$ws.Range("M122").Value = -4323.000100000001
$ws.Range("H136").Value = 798.75
$ws.Range("I136").Value = 798.75
$ws.Range("K136").Value = 2396.25
$ws.Range("M136").Value = 153.75

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 29999
$ws.Range("J125").Value = 29999
$ws.Range("L125").Value = 29999
$ws.Range("N125").Value = -39839

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1519.1333
$ws.Range("I31").Value = 1484.7858
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1484.7858
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -1189.7858
$ws.Range("N31").Value = -2590
$ws.Range("H34").Value = 1519.1333
$ws.Range("I34").Value = 1484.7858
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 1484.7858
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -1282.7858
$ws.Range("N34").Value = -2404
$ws.Range("H41").Value = 6581
$ws.Range("I41").Value = 3224.5
$ws.Range("J41").Value = 9937.5
$ws.Range("K41").Value = 3224.5
$ws.Range("L41").Value = 9937.5
$ws.Range("M41").Value = -2796.5
$ws.Range("N41").Value = -10793.5
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0
$ws.Range("H132").Value = 2198.6072
$ws.Range("I132").Value = 2261.6296
$ws.Range("J132").Value = 497
$ws.Range("K132").Value = 6784.888800000001
$ws.Range("L132").Value = 1491
$ws.Range("M132").Value = -4254.888800000001
$ws.Range("N132").Value = -6551
$ws.Range("H134").Value = 2189.8
$ws.Range("I134").Value = 1621.1428
$ws.Range("J134").Value = 3516.6667
$ws.Range("K134").Value = 4863.428400000001
$ws.Range("L134").Value = 10550.0001
$ws.Range("M134").Value = -2328.428400000001
$ws.Range("N134").Value = -15620.0001

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9668.666999999999
$ws.Range("I80").Value = 8500
$ws.Range("K80").Value = 8500
$ws.Range("M80").Value = -7502
$ws.Range("H83").Value = 9668.666999999999
$ws.Range("I83").Value = 8500
$ws.Range("K83").Value = 42500
$ws.Range("M83").Value = -37508
$ws.Range("H107").Value = 6093.9
$ws.Range("I107").Value = 1166.3334
$ws.Range("J107").Value = 13485.25
$ws.Range("K107").Value = 1166.3334
$ws.Range("L107").Value = 13485.25
$ws.Range("M107").Value = 753.6666
$ws.Range("N107").Value = -17325.25
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1544.9166
$ws.Range("I93").Value = 1130.625
$ws.Range("J93").Value = 2373.5
$ws.Range("K93").Value = 1130.625
$ws.Range("L93").Value = 2373.5
$ws.Range("M93").Value = 117.375
$ws.Range("N93").Value = -4869.5
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H132").Value = 1686.6
$ws.Range("I132").Value = 1686.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5059.799999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -2529.799999999999

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1880
$ws.Range("I96").Value = 1880
$ws.Range("K96").Value = 1880
$ws.Range("M96").Value = -507
$ws.Range("H103").Value = 11900.25
$ws.Range("J103").Value = 11900.25
$ws.Range("L103").Value = 11900.25
$ws.Range("N103").Value = -14244.25
$ws.Range("H107").Value = 609.8889
$ws.Range("J107").Value = 745
$ws.Range("L107").Value = 2235
$ws.Range("N107").Value = -6075
$ws.Range("H132").Value = 1089.2307
$ws.Range("I132").Value = 1121.9445
$ws.Range("K132").Value = 3365.8335
$ws.Range("M132").Value = -835.8335000000002
